$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'68.281.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "'  -0.16%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value2 = "'2.645.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "'  +0.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value2 = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value2 = "'597.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "'  -0.33%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value2 = "'157.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "'  +2.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value2 = "'0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "'  -0.65%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value2 = "'0.142"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "'  +2.63%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value2 = "'  -1.47%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value2 = "'5.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "'  +0.61%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value2 = "'  +0.83%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value2 = "'28.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "'  +0.41%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value2 = "'  +0.85%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value2 = "'3.125.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value2 = "'68.147.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "'  -0.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value2 = "'2.629.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "'  -0.37%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value2 = "'11.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "'  -0.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value2 = "'364.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "'  -0.45%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value2 = "'Polkadot"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value2 = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value2 = "'4.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "'  +3.80%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value2 = "'Uniswap"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value2 = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value2 = "'7.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "'  -1.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value2 = "'4.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "'  -0.89%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value2 = "'  -1.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value2 = "'75.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "'  +2.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value2 = "'  +0.07%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value2 = "'9.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "'  -2.72%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value2 = "'2.784.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "'  +0.42%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value2 = "'  -0.18%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value2 = "'  +0.13%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value2 = "'557.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "'  -3.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value2 = "'8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "'  +0.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value2 = "'  -0.89%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E34").Value2 = "'  -0.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value2 = "'  -0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value2 = "'1.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "'  +1.44%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value2 = "'  +2.90%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value2 = "'159.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "'  -0.29%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value2 = "'  +0.71%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value2 = "'  -2.31%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value2 = "'5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "'  -0.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value2 = "'  +3.96%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value2 = "'2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "'  -0.55%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value2 = "'  +0.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value2 = "'158.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "'  +0.90%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value2 = "'3.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "'  +0.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value2 = "'22.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "'  +1.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value2 = "'1.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "'  -1.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value2 = "'  +0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value2 = "'0.615"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "'  +0.01%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value2 = "'0.566"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "'  +0.90%  "
$ws.Range("E51").Style = "Normal"
